# Update "想去人数" (interested-people count) values in the
# "展览" and "全部类型" worksheets to reflect refreshed counts.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (row -> new value for column F)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3975
$ws1.Range("F4").Value = 2335
$ws1.Range("F7").Value = 30
$ws1.Range("F8").Value = 14
$ws1.Range("F11").Value = 48
$ws1.Range("F13").Value = 1479
$ws1.Range("F14").Value = 262
$ws1.Range("F15").Value = 2738
$ws1.Range("F16").Value = 189

# Sheet "全部类型" (row -> new value for column F)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3975
$ws4.Range("F4").Value = 2335
$ws4.Range("F7").Value = 30
$ws4.Range("F8").Value = 14
$ws4.Range("F12").Value = 48
$ws4.Range("F16").Value = 1479
$ws4.Range("F17").Value = 262
$ws4.Range("F18").Value = 2738
$ws4.Range("F19").Value = 189
